$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: date label (col A, stored as text) + values for columns B..J
$data = @(
    @("01-07-2021", 191, 0, 0, 0, 27, 39, 68, 23, 35),
    @("02-07-2021", 169, 0, 0, 0, 0, 58, 37, 58, 16),
    @("05-07-2021", 3, 0, 0, 0, 0, 0, 0, 0, 3),
    @("06-07-2021", 342, 0, 168, 0, 0, 0, 80, 68, 26),
    @("07-07-2021", 83, 0, 0, 0, 0, 0, 32, 11, 40),
    @("08-07-2021", 338, 192, 0, 30, 45, 5, 31, 18, 17),
    @("09-07-2021", 128, 0, 48, 0, 0, 25, 25, 13, 17),
    @("12-07-2021", 415, 0, 0, 64, 121, 0, 138, 87, 6),
    @("13-07-2021", 532, 192, 48, 32, 44, 15, 162, 33, 7),
    @("14-07-2021", 479, 89, 15, 0, 202, 16, 67, 49, 42),
    @("15-07-2021", 525, 294, 48, 0, 22, 0, 16, 75, 69),
    @("19-07-2021", 72, 0, 0, 0, 38, 0, 8, 10, 16),
    @("20-07-2021", 305, 0, 74, 0, 200, 0, 13, 12, 6),
    @("21-07-2021", 85, 0, 0, 0, 23, 0, 38, 16, 8),
    @("22-07-2021", 140, 0, 0, 60, 38, 0, 38, 0, 4),
    @("23-07-2021", 69, 0, 0, 0, 15, 0, 38, 5, 11),
    @("26-07-2021", 185, 0, 0, 0, 38, 81, 44, 5, 16),
    @("27-07-2021", 295, 0, 96, 32, 38, 51, 19, 25, 35),
    @("28-07-2021", 424, 96, 0, 0, 124, 25, 146, 11, 21),
    @("29-07-2021", 169, 0, 0, 26, 48, 0, 69, 3, 23),
    @("30-07-2021", 443, 288, 48, 50, 22, 0, 12, 4, 19)
)

$startRow = 126

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    # Column A holds a date-looking label that must stay plain text (not get
    # auto-converted into a date serial by the "smart" value setter). Stage it
    # in a scratch cell as a text formula, then copy/paste-values it into
    # place -- this preserves the literal string with no number-format churn.
    $ws.Range("Z1").Formula = "=""" + $rowData[0] + """"
    $ws.Range("Z1").Copy()
    $ws.Range("A" + $row).PasteSpecial(-4163)

    for ($c = 1; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}

$ws.Range("Z1").Value = ""
